# Update absenteeism data rows 2-11 with new values per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 67121
$ws.Range("B2").Value = "Bruno Martins"
$ws.Range("C2").Value = "Marketing"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45087
$ws.Range("G2").Value = 11583.25

# Row 3
$ws.Range("A3").Value = 18464
$ws.Range("B3").Value = "Bárbara Moura"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45102
$ws.Range("G3").Value = 10236.32

# Row 4
$ws.Range("A4").Value = 85612
$ws.Range("B4").Value = "Gustavo Viana"
$ws.Range("C4").Value = "Operações"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 5341.72

# Row 5
$ws.Range("A5").Value = 32708
$ws.Range("B5").Value = "André Rodrigues"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45092
$ws.Range("G5").Value = 11480.36

# Row 6
$ws.Range("A6").Value = 13728
$ws.Range("B6").Value = "Davi Rocha"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45085
$ws.Range("G6").Value = 10556.85

# Row 7
$ws.Range("A7").Value = 9102
$ws.Range("B7").Value = "Joana da Mata"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 7242.48

# Row 8
$ws.Range("A8").Value = 52940
$ws.Range("B8").Value = "Beatriz Gonçalves"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 9763.23

# Row 9
$ws.Range("A9").Value = 43285
$ws.Range("B9").Value = "Dr. Daniel Alves"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 2823.34

# Row 10
$ws.Range("A10").Value = 59511
$ws.Range("B10").Value = "Esther Azevedo"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 9750.9

# Row 11
$ws.Range("A11").Value = 65143
$ws.Range("B11").Value = "Bernardo Fogaça"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45098
$ws.Range("G11").Value = 6520.38
